# Apply updated crypto price / volume figures to columns D and E
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.783.37'
$ws.Range('E2').Value = '  +5.47%  '
$ws.Range('D3').Value = '1.705.38'
$ws.Range('E3').Value = '  +3.37%  '
$ws.Range('D4').Value = '''1.002'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''331.22'
$ws.Range('E5').Value = '  +6.25%  '
$ws.Range('D6').Value = '''1.000'
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('D7').Value = '''0.3685'
$ws.Range('E7').Value = '  +0.88%  '
$ws.Range('D8').Value = '''48.50'
$ws.Range('E8').Value = '  +4.33%  '
$ws.Range('D9').Value = '''0.3305'
$ws.Range('E9').Value = '  +1.87%  '
$ws.Range('D10').Value = '''1.169'
$ws.Range('E10').Value = '  +4.11%  '
$ws.Range('D11').Value = '''0.07334'
$ws.Range('E11').Value = '  +4.56%  '
$ws.Range('D12').Value = '''1.000'
$ws.Range('E12').Value = '  +0.08%  '
$ws.Range('D13').Value = '''6.194'
$ws.Range('E13').Value = '  +3.88%  '
$ws.Range('D14').Value = '''20.01'
$ws.Range('E14').Value = '  +3.53%  '
$ws.Range('D15').Value = '''6.859'
$ws.Range('E15').Value = '  +3.86%  '
$ws.Range('D16').Value = '1.702.98'
$ws.Range('E16').Value = '  +2.95%  '
$ws.Range('D17').Value = '''0.00001064'
$ws.Range('E17').Value = '  +2.44%  '
$ws.Range('D18').Value = '''0.06634'
$ws.Range('D19').Value = '''81.14'
$ws.Range('E19').Value = '  +3.46%  '
$ws.Range('D20').Value = '''0.9998'
$ws.Range('E20').Value = '  +0.13%  '
$ws.Range('D21').Value = '''16.15'
$ws.Range('E21').Value = '  +3.34%  '
$ws.Range('D22').Value = '''6.045'
$ws.Range('E22').Value = '  +2.15%  '
$ws.Range('D23').Value = '''12.99'
$ws.Range('E23').Value = '  +3.65%  '
$ws.Range('D24').Value = '25.741.71'
$ws.Range('E24').Value = '  +5.40%  '
$ws.Range('D25').Value = '''2.457'
$ws.Range('E25').Value = '  -0.11%  '
$ws.Range('D26').Value = '''2.476'
$ws.Range('E26').Value = '  +6.82%  '
$ws.Range('D27').Value = '''149.69'
$ws.Range('E27').Value = '  +2.07%  '
$ws.Range('D28').Value = '''19.17'
$ws.Range('E28').Value = '  +3.27%  '
$ws.Range('D29').Value = '''1.289'
$ws.Range('E29').Value = '  +8.28%  '
$ws.Range('D30').Value = '1.890.30'
$ws.Range('E30').Value = '  +3.24%  '
$ws.Range('D31').Value = '''128.06'
$ws.Range('E31').Value = '  +3.26%  '
$ws.Range('D32').Value = '''4.090'
$ws.Range('E32').Value = '  +0.79%  '
$ws.Range('D33').Value = '''5.928'
$ws.Range('E33').Value = '  +4.11%  '
$ws.Range('D34').Value = '''1.718'
$ws.Range('E34').Value = '  +3.63%  '
$ws.Range('D35').Value = '''0.08491'
$ws.Range('E35').Value = '  +0.57%  '
$ws.Range('D36').Value = '''12.85'
$ws.Range('E36').Value = '  +6.33%  '
$ws.Range('D37').Value = '''5.328'
$ws.Range('E37').Value = '  +2.45%  '
$ws.Range('D38').Value = '''1.272'
$ws.Range('E38').Value = '  +0.48%  '
$ws.Range('D39').Value = '''0.06190'
$ws.Range('E39').Value = '  +2.87%  '
$ws.Range('D40').Value = '''8.521'
$ws.Range('E40').Value = '  +5.45%  '
$ws.Range('D41').Value = '''0.2114'
$ws.Range('E41').Value = '  +2.53%  '
$ws.Range('D42').Value = '''0.02248'
$ws.Range('E42').Value = '  +0.71%  '
$ws.Range('D43').Value = '''14.58'
$ws.Range('E43').Value = '  +16.04%  '
$ws.Range('D44').Value = '''0.6109'
$ws.Range('E44').Value = '  +3.76%  '
$ws.Range('D45').Value = '''1.000'
$ws.Range('D46').Value = '''3.844'
$ws.Range('E46').Value = '  +2.17%  '
$ws.Range('D47').Value = '''0.5823'
$ws.Range('E47').Value = '  +3.84%  '
$ws.Range('D48').Value = '''126.66'
$ws.Range('E48').Value = '  +3.08%  '
$ws.Range('D49').Value = '''1.999'
$ws.Range('E49').Value = '  +2.88%  '
$ws.Range('D50').Value = '''0.07214'
$ws.Range('E50').Value = '  +4.53%  '
$ws.Range('D51').Value = '''1.205'
$ws.Range('E51').Value = '  +2.14%  '

Write-Host "Updated cryptos list"
